$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 5941.125
$ws.Range("I86").Value = 2923.375
$ws.Range("J86").Value = 8958.875
$ws.Range("K86").Value = 2923.375
$ws.Range("L86").Value = 8958.875
$ws.Range("M86").Value = -1800.375
$ws.Range("N86").Value = -11204.875

$ws.Range("H89").Value = 5941.125
$ws.Range("I89").Value = 2923.375
$ws.Range("J89").Value = 8958.875
$ws.Range("K89").Value = 14616.875
$ws.Range("L89").Value = 44794.375
$ws.Range("M89").Value = -9000.875
$ws.Range("N89").Value = -56026.375

$ws.Range("H98").Value = 1658.742
$ws.Range("I98").Value = 984.625
$ws.Range("J98").Value = 3970
$ws.Range("K98").Value = 984.625
$ws.Range("L98").Value = 3970
$ws.Range("M98").Value = 513.375
$ws.Range("N98").Value = -6966

$ws.Range("H122").Value = 1658.742
$ws.Range("I122").Value = 984.625
$ws.Range("J122").Value = 3970
$ws.Range("K122").Value = 2953.875
$ws.Range("L122").Value = 11910
$ws.Range("M122").Value = -503.875
$ws.Range("N122").Value = -16810

$ws.Range("H132").Value = 2221.228
$ws.Range("I132").Value = 1304.5333
$ws.Range("J132").Value = 5658.8335
$ws.Range("K132").Value = 3913.5999
$ws.Range("L132").Value = 16976.5005
$ws.Range("M132").Value = -1383.5999
$ws.Range("N132").Value = -22036.5005

$ws.Range("H133").Value = 31791.666
$ws.Range("J133").Value = 31791.666
$ws.Range("L133").Value = 31791.666
$ws.Range("N133").Value = -41911.666

$ws.Range("H137").Value = 3610.3225
$ws.Range("I137").Value = 3639.8696
$ws.Range("J137").Value = 3525.375
$ws.Range("K137").Value = 10919.6088
$ws.Range("L137").Value = 10576.125
$ws.Range("M137").Value = -8369.6088
$ws.Range("N137").Value = -15676.125

$ws.Range("H141").Value = 2577.9614
$ws.Range("I141").Value = 1031.9474
$ws.Range("K141").Value = 3095.8422
$ws.Range("M141").Value = 2084.1578

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 19229.795
$ws.Range("I32").Value = 18937.623
$ws.Range("J32").Value = 21469.777
$ws.Range("K32").Value = 18937.623
$ws.Range("L32").Value = 21469.777
$ws.Range("M32").Value = -18650.623
$ws.Range("N32").Value = -22043.777

$ws.Range("H97").Value = 673.5
$ws.Range("I97").Value = 704.375
$ws.Range("J97").Value = 550
$ws.Range("K97").Value = 704.375
$ws.Range("L97").Value = 550
$ws.Range("M97").Value = -208.375
$ws.Range("N97").Value = -1542

$ws.Range("H122").Value = 3995.5144
$ws.Range("I122").Value = 4208.409
$ws.Range("J122").Value = 3635.2307
$ws.Range("K122").Value = 12625.227
$ws.Range("L122").Value = 10905.6921
$ws.Range("M122").Value = -10175.227
$ws.Range("N122").Value = -15805.6921

$ws.Range("H132").Value = 18187.203
$ws.Range("I132").Value = 22923.438
$ws.Range("J132").Value = 3978.5
$ws.Range("K132").Value = 68770.314
$ws.Range("L132").Value = 11935.5
$ws.Range("M132").Value = -66240.314
$ws.Range("N132").Value = -16995.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 4531.875
$ws.Range("I99").Value = 4843.5713
$ws.Range("J99").Value = 2350
$ws.Range("K99").Value = 4843.5713
$ws.Range("L99").Value = 2350
$ws.Range("M99").Value = -3345.5713
$ws.Range("N99").Value = -5346

$ws.Range("H134").Value = 2376.2537
$ws.Range("I134").Value = 1944.9814
$ws.Range("J134").Value = 4167.6924
$ws.Range("K134").Value = 5834.9442
$ws.Range("L134").Value = 12503.0772
$ws.Range("M134").Value = -3299.9442
$ws.Range("N134").Value = -17573.0772

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2545.2173
$ws.Range("I31").Value = 1882.7678
$ws.Range("J31").Value = 5398.846
$ws.Range("K31").Value = 1882.7678
$ws.Range("L31").Value = 5398.846
$ws.Range("M31").Value = -1587.7678
$ws.Range("N31").Value = -5988.846

$ws.Range("H34").Value = 2545.2173
$ws.Range("I34").Value = 1882.7678
$ws.Range("J34").Value = 5398.846
$ws.Range("K34").Value = 1882.7678
$ws.Range("L34").Value = 5398.846
$ws.Range("M34").Value = -1680.7678
$ws.Range("N34").Value = -5802.846

$ws.Range("H122").Value = 1581.421
$ws.Range("I122").Value = 959.3333
$ws.Range("K122").Value = 2877.9999
$ws.Range("M122").Value = -427.9998999999998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 3275.5557
$ws.Range("I63").Value = 870
$ws.Range("J63").Value = 5200
$ws.Range("K63").Value = 2610
$ws.Range("L63").Value = 15600
$ws.Range("M63").Value = -1861
$ws.Range("N63").Value = -17098

$ws.Range("H66").Value = 3275.5557
$ws.Range("I66").Value = 870
$ws.Range("J66").Value = 5200
$ws.Range("K66").Value = 7830
$ws.Range("L66").Value = 46800
$ws.Range("M66").Value = -4086
$ws.Range("N66").Value = -54288

$ws.Range("H113").Value = 546.9153
$ws.Range("I113").Value = 560.6
$ws.Range("J113").Value = 502.92856
$ws.Range("K113").Value = 1681.8
$ws.Range("L113").Value = 1508.78568
$ws.Range("M113").Value = 488.1999999999998
$ws.Range("N113").Value = -5848.78568

$ws.Range("H131").Value = 1604.9762
$ws.Range("I131").Value = 1714.125
$ws.Range("J131").Value = 1579.2941
$ws.Range("K131").Value = 5142.375
$ws.Range("L131").Value = 4737.8823
$ws.Range("M131").Value = -102.375
$ws.Range("N131").Value = -14817.8823

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 334567.56
$ws.Range("I122").Value = 429472.56
$ws.Range("J122").Value = 2400
$ws.Range("K122").Value = 1288417.68
$ws.Range("L122").Value = 7200
$ws.Range("M122").Value = -1285967.68
$ws.Range("N122").Value = -12100

$ws.Range("H132").Value = 3511.4905
$ws.Range("I132").Value = 3413.1428
$ws.Range("J132").Value = 3702.7222
$ws.Range("K132").Value = 10239.4284
$ws.Range("L132").Value = 11108.1666
$ws.Range("M132").Value = -7709.428400000001
$ws.Range("N132").Value = -16168.1666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2888.3462
$ws.Range("I7").Value = 2883.2
$ws.Range("J7").Value = 2895.3635
$ws.Range("K7").Value = 2883.2
$ws.Range("L7").Value = 2895.3635
$ws.Range("M7").Value = -2771.2
$ws.Range("N7").Value = -3119.3635

$ws.Range("H122").Value = 1962.5
$ws.Range("I122").Value = 1814.2858
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 5442.857400000001
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -2992.857400000001
$ws.Range("N122").Value = -13900

$ws.Range("H126").Value = 2888.3462
$ws.Range("I126").Value = 2883.2
$ws.Range("J126").Value = 2895.3635
$ws.Range("K126").Value = 8649.599999999999
$ws.Range("L126").Value = 8686.0905
$ws.Range("M126").Value = -6179.599999999999
$ws.Range("N126").Value = -13626.0905

$ws.Range("H134").Value = 24684.75
$ws.Range("J134").Value = 24684.75
$ws.Range("L134").Value = 24684.75
$ws.Range("N134").Value = -34824.75

$ws.Range("H137").Value = 29901.334
$ws.Range("J137").Value = 31481.8
$ws.Range("L137").Value = 31481.8
$ws.Range("N137").Value = -41681.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 33367976
$ws.Range("I122").Value = 43522984
$ws.Range("J122").Value = 1515
$ws.Range("K122").Value = 130568952
$ws.Range("L122").Value = 4545
$ws.Range("M122").Value = -130566502
$ws.Range("N122").Value = -9445
